$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet
$ws.Name = "overlap with PMID- 33444546"

# New column header (copy the header style from E1 so F1 matches the other headers)
$ws.Range("F1").Value = "well_determined_region_names"
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# New column F data per row
$regions = @(
    "Region 2 (Nsp2; 74.77%)",
    "Region 16 (Nsp3; 100.00%)",
    "Region 25 (Nsp8, Nsp9; 65.46%)",
    "Region 25 (Nsp8, Nsp9; 100.00%)",
    "None",
    "None",
    "Region 36 (Nsp13; 75.22%)",
    "Region 37 (Nsp14; 100.00%)",
    "N/A",
    "N/A",
    "N/A"
)

for ($i = 0; $i -lt $regions.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $regions[$i]
}

